$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header and municipality/state name values (title-case "De/Del/El/La/Los/Y")
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B4').Value = 'Rincón De Romos'
$ws.Range('B16').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Mazapa De Madero'
$ws.Range('B60').Value = 'San Juan De Sabinas'
$ws.Range('A67').Value = 'Ciudad De México'
$ws.Range('B71').Value = 'Cuajimalpa De Morelos'
$ws.Range('A91').Value = 'Estado De México'
$ws.Range('B92').Value = 'Almoloya De Juárez'
$ws.Range('B95').Value = 'Atizapán De Zaragoza'
$ws.Range('B101').Value = 'Ecatepec De Morelos'
$ws.Range('B104').Value = 'Ixtapan De La Sal'
$ws.Range('B108').Value = 'Naucalpan De Juárez'
$ws.Range('B110').Value = 'San Felipe Del Progreso'
$ws.Range('B116').Value = 'Tenango Del Valle'
$ws.Range('B120').Value = 'Tlalnepantla De Baz'
$ws.Range('B124').Value = 'Villa Del Carbón'
$ws.Range('B130').Value = 'Apaseo El Alto'
$ws.Range('B131').Value = 'Apaseo El Grande'
$ws.Range('B136').Value = 'Jaral Del Progreso'
$ws.Range('B140').Value = 'Purísima Del Rincón'
$ws.Range('B143').Value = 'San Diego De La Unión'
$ws.Range('B146').Value = 'San Luis De La Paz'
$ws.Range('B147').Value = 'Valle De Santiago'
$ws.Range('B150').Value = 'Acapulco De Juárez'
$ws.Range('B152').Value = 'Alcozauca De Guerrero'
$ws.Range('B155').Value = 'Ayutla De Los Libres'
$ws.Range('B158').Value = 'Chilapa De Álvarez'
$ws.Range('B159').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B160').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B161').Value = 'Coyuca De Benítez'
$ws.Range('B162').Value = 'Coyuca De Catalán'
$ws.Range('B164').Value = 'Cutzamala De Pinzón'
$ws.Range('B167').Value = 'Iguala De La Independencia'
$ws.Range('B168').Value = 'Zihuatanejo De Azueta'
$ws.Range('B177').Value = 'Taxco De Alarcón'
$ws.Range('B179').Value = 'Técpan De Galeana'
$ws.Range('B181').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B184').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B185').Value = 'Tlapa De Comonfort'
$ws.Range('B195').Value = 'Atotonilco El Grande'
$ws.Range('B202').Value = 'Jacala De Ledezma'
$ws.Range('B208').Value = 'Omitlán De Juárez'
$ws.Range('B209').Value = 'Pachuca De Soto'
$ws.Range('B214').Value = 'Tepehuacán De Guerrero'
$ws.Range('B215').Value = 'Tezontepec De Aldama'
$ws.Range('B217').Value = 'Tula De Allende'
$ws.Range('B218').Value = 'Tulancingo De Bravo'
$ws.Range('B220').Value = 'Zacualtipán De Ángeles'
$ws.Range('B221').Value = 'Zapotlán De Juárez'
$ws.Range('B224').Value = 'Acatlán De Juárez'
$ws.Range('B228').Value = 'Encarnación De Díaz'
$ws.Range('B230').Value = 'Lagos De Moreno'
$ws.Range('B232').Value = 'San Juan De Los Lagos'
$ws.Range('B233').Value = 'Tamazula De Gordiano'
$ws.Range('B235').Value = 'Tepatitlán De Morelos'
$ws.Range('B236').Value = 'Tizapán El Alto'
$ws.Range('B238').Value = 'Unión De San Antonio'
$ws.Range('B269').Value = 'Tetela Del Volcán'
$ws.Range('B284').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B289').Value = 'Fresnillo De Trujano'
$ws.Range('B290').Value = 'Guevea De Humboldt'
$ws.Range('B291').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B292').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B293').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B294').Value = 'Ixtlán De Juárez'
$ws.Range('B295').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B297').Value = 'Mariscala De Juárez'
$ws.Range('B298').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B299').Value = 'Oaxaca De Juárez'
$ws.Range('B300').Value = 'Ocotlán De Morelos'
$ws.Range('B301').Value = 'Pinotepa De Don Luis'
$ws.Range('B302').Value = 'Putla Villa De Guerrero'
$ws.Range('B312').Value = 'San Francisco Del Mar'
$ws.Range('B314').Value = 'San José Del Progreso'
$ws.Range('B320').Value = 'San Juan Del Río'
$ws.Range('B343').Value = 'Santa Ana Del Valle'
$ws.Range('B364').Value = 'Teotitlán Del Valle'
$ws.Range('B365').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B366').Value = 'Tlacolula De Matamoros'
$ws.Range('B367').Value = 'Totontepec Villa De Morelos'
$ws.Range('B368').Value = 'Villa De Etla'
$ws.Range('B369').Value = 'Villa De Zaachila'
$ws.Range('B370').Value = 'Villa Sola De Vega'
$ws.Range('B378').Value = 'Ayotoxco De Guerrero'
$ws.Range('B387').Value = 'Cuayuca De Andrade'
$ws.Range('B394').Value = 'Huehuetlán El Grande'
$ws.Range('B396').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B398').Value = 'Izúcar De Matamoros'
$ws.Range('B402').Value = 'Los Reyes De Juárez'
$ws.Range('B405').Value = 'Palmar De Bravo'
$ws.Range('B411').Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Range('B420').Value = 'San Salvador El Seco'
$ws.Range('B426').Value = 'Tepanco De López'
$ws.Range('B427').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B432').Value = 'Tetela De Ocampo'
$ws.Range('B435').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B452').Value = 'Jalpan De Serra'
$ws.Range('B453').Value = 'Landa De Matamoros'
$ws.Range('B455').Value = 'San Juan Del Río'
$ws.Range('B462').Value = 'Ciudad Del Maíz'
$ws.Range('B468').Value = 'San Ciro De Acosta'
$ws.Range('B474').Value = 'Villa De La Paz'
$ws.Range('B475').Value = 'Villa De Ramos'
$ws.Range('B518').Value = 'Boca Del Río'
$ws.Range('B520').Value = 'Cazones De Herrera'
$ws.Range('B524').Value = 'Cosamaloapan De Carpio'
$ws.Range('B530').Value = 'Huiloapan De Cuauhtémoc'
$ws.Range('B531').Value = 'Ignacio De La Llave'
$ws.Range('B533').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B534').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B537').Value = 'Juchique De Ferrer'
$ws.Range('B541').Value = 'Martínez De La Torre'
$ws.Range('B554').Value = 'Poza Rica De Hidalgo'
$ws.Range('B564').Value = 'Tlacotepec De Mejía'
$ws.Range('B577').Value = 'Noria De Ángeles'

# Delete the trailing metadata rows (582:586) and update the used range/dimension
$ws.Range("A582:A586").EntireRow.Delete() | Out-Null

